$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Set-NumCellStyled($cellRef, $val, $styleDonor) {
    $ws.Range($cellRef).Value = $val
    $ws.Range($styleDonor).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

function Set-TextCell($cellRef, $text, $styleDonor) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($styleDonor).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- Header text updates (Volume/Number + report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"

# --- Crime Complaints table updates (rows 15-33) ---

Set-TextCell "D15" "0" "A14"
Set-TextCell "E15" "***.*" "A14"
Set-NumCell "F15" 2
Set-NumCell "H15" 0
Set-NumCell "N15" -33.333333333333
Set-TextCell "D16" "0" "A14"
Set-TextCell "E16" "***.*" "A14"
Set-NumCell "F16" 5
Set-NumCell "G16" 3
Set-NumCell "H16" 66.666666666666
Set-NumCell "I16" 46
Set-NumCell "K16" -13.207547169811
Set-NumCell "L16" 0
Set-NumCell "M16" -55.339805825242
Set-NumCell "N16" -83.333333333333
Set-NumCell "C17" 6
Set-NumCell "D17" 6
Set-NumCell "E17" 0
Set-NumCell "G17" 12
Set-NumCell "H17" 83.333333333333
Set-NumCell "I17" 233
Set-NumCell "J17" 147
Set-NumCell "K17" 58.503401360544
Set-NumCell "L17" 52.287581699346
Set-NumCell "M17" 86.4
Set-NumCell "N17" -17.957746478873
Set-NumCell "C18" 3
Set-NumCell "E18" 50
Set-NumCell "F18" 7
Set-NumCell "G18" 9
Set-NumCell "H18" -22.222222222222
Set-NumCell "I18" 83
Set-NumCell "J18" 93
Set-NumCell "K18" -10.752688172043
Set-NumCell "L18" -1.190476190476
Set-NumCell "M18" -55.851063829787
Set-NumCell "N18" -93.241042345276
Set-NumCell "C19" 14
Set-NumCell "D19" 8
Set-NumCell "E19" 75
Set-NumCell "F19" 45
Set-NumCell "G19" 29
Set-NumCell "H19" 55.172413793103
Set-NumCell "I19" 375
Set-NumCell "J19" 427
Set-NumCell "K19" -12.177985948477
Set-NumCell "L19" -13.594470046083
Set-NumCell "M19" 0.53619302949
Set-NumCell "N19" -50.331125827814
Set-NumCell "F20" 6
Set-NumCell "H20" 50
Set-NumCell "I20" 51
Set-NumCell "J20" 72
Set-NumCell "K20" -29.166666666666
Set-NumCell "L20" -48.484848484848
Set-NumCell "M20" -51.428571428571
Set-NumCell "N20" -97.974583002382
Set-NumCell "C21" 27
Set-NumCell "D21" 17
Set-NumCell "E21" 58.823529411764
Set-NumCell "F21" 87
Set-NumCell "G21" 59
Set-NumCell "H21" 47.457627118644
Set-NumCell "I21" 800
Set-NumCell "J21" 809
Set-NumCell "K21" -1.112484548825
Set-NumCell "L21" -3.147699757869
Set-NumCell "M21" -12.280701754386
Set-NumCell "N21" -84.258166076347
Set-NumCellStyled "C23" 1 "G15"
Set-NumCell "D23" 2
Set-NumCell "E23" -50
Set-NumCell "G23" 4
Set-NumCell "H23" -50
Set-NumCell "I23" 55
Set-NumCell "J23" 26
Set-NumCell "K23" 111.538461538462
Set-NumCell "L23" 71.875
Set-NumCell "M23" 161.904761904762
Set-NumCell "C24" 19
Set-NumCell "D24" 15
Set-NumCell "E24" 26.666666666666
Set-NumCell "F24" 68
Set-NumCell "G24" 80
Set-NumCell "H24" -15
Set-NumCell "I24" 914
Set-NumCell "J24" 907
Set-NumCell "K24" 0.77177508269
Set-NumCell "L24" -6.639427987742
Set-NumCell "M24" -40.417209908735
Set-NumCell "D25" 6
Set-NumCell "E25" 0
Set-NumCell "F25" 31
Set-NumCell "G25" 49
Set-NumCell "H25" -36.734693877551
Set-NumCell "I25" 470
Set-NumCell "J25" 442
Set-NumCell "K25" 6.334841628959
Set-NumCell "L25" 17.5
Set-NumCell "C26" 4
Set-NumCell "E26" -20
Set-NumCell "F26" 21
Set-NumCell "G26" 21
Set-NumCell "H26" 0
Set-NumCell "I26" 326
Set-NumCell "J26" 327
Set-NumCell "K26" -0.305810397553
Set-NumCell "L26" 8.666666666666
Set-NumCell "M26" -33.05954825462
Set-TextCell "C27" "0" "A14"
Set-TextCell "D27" "0" "A14"
Set-TextCell "E27" "***.*" "A14"
Set-NumCell "F27" 3
Set-NumCell "H27" 50
Set-TextCell "C28" "0" "A14"
Set-NumCell "F28" 4
Set-NumCell "G28" 1
Set-NumCell "H28" 300
Set-NumCell "I28" 38
Set-NumCell "K28" 11.764705882352
Set-NumCell "L28" 8.571428571428
Set-TextCell "F31" "0" "A14"
Set-TextCell "D33" "0" "A14"
Set-TextCell "E33" "***.*" "A14"

$excel.CutCopyMode = 0
Write-Output "done"
